# Added 1.1.0 of term
#
# 1. Bump the "Version" metadata row from 1.0.0 -> 1.1.0
# 2. Bump the "Date" metadata row to the new publish timestamp
# 3. Make the header/body cell styles actually *apply* their alignment
#    (wrapText + vertical-top) by turning WrapText / VerticalAlignment on
#    for the header row and the body rows of every sheet - this is what
#    flips the corresponding cellXfs records to applyAlignment="true".

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")

# "Version" is row 3 (A3 label / B3 value)
$metadata.Range("B3").Value = "1.1.0"

# "Date" is row 8 (A8 label / B8 value)
$metadata.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Re-assert the (previously inert) wrap/vertical-top alignment on the
# header row and the body of every worksheet so it is actually applied.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    $headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $colCount))
    $headerRange.VerticalAlignment = -4160
    $headerRange.WrapText = $true

    if ($rowCount -gt 1) {
        $bodyRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($rowCount, $colCount))
        $bodyRange.VerticalAlignment = -4160
        $bodyRange.WrapText = $true
    }
}
